# "DeudoresPrueba" refresh: replace the debtor list with the latest export.
# The new export has 27 rows (vs. the previous 38), with several clients
# renamed/removed and all Fecha/Valor figures updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing rows that no longer exist in the refreshed export
# (old data ran through row 39, new data stops at row 28).
$ws.Range("A29:E39").EntireRow.Delete()

# New row data: Consecutivo (A) is left as the existing 1..27 sequence,
# Pagado (E) stays FALSE for every row, so only Cliente/Fecha/Valor change.
$rows = @(
    @(2, "ALISO", 46032, 147000),
    @(3, "ARROZ PAISA SUBA", 46029, 166000),
    @(4, "CAMPO VERDE ZIPAQUIRA", 46021, 18900),
    @(5, "CARNES JOHANA", 46035, 164000),
    @(6, "CARNILANDIA", 46028, 100000),
    @(7, "CIMARRON DORADO", 46035, 797000),
    @(8, "COCINA CHINA", 46031, 170000),
    @(9, "COCINA CHINA", 46036, 170000),
    @(10, "DARWIN FUTBOL", 45921, 200000),
    @(11, "DAVIDCITO", 45947, 100000),
    @(12, "EL JORDAN", 46022, 1600000),
    @(13, "FRANCO", 45996, 20000),
    @(14, "FRANCO", 46017, 545800),
    @(15, "LA SELECTA", 45912, 82000),
    @(16, "LINARIA", 46034, 185000),
    @(17, "MAFE", 46034, 521300),
    @(18, "MERKA FRUVER DEXI", 45995, 339000),
    @(19, "MERKA FRUVER DEXI", 45988, 15400),
    @(20, "NEVADA", 46031, 321900),
    @(21, "PARAISO FUNZA", 46035, 100000),
    @(22, "PINILLA", 45931, 82000),
    @(23, "PLAZA JESSICA", 46028, 285500),
    @(24, "PLAZA JESSICA", 46033, 702500),
    @(25, "PUNTA DE ANCA", 46031, 257000),
    @(26, "EL RUBY", 46035, 110600),
    @(27, "CALDAS WOK", 46035, 85000),
    @(28, "CARNIVOROS", 46036, 196800)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
}

# Cosmetic touches matching the refreshed workbook: best-fit column widths
# for Cliente/Fecha and the selection left over from editing.
$ws.Columns.Item(2).ColumnWidth = 25.7109375
$ws.Columns.Item(3).ColumnWidth = 10.42578125
$null = $ws.Range("A2:A28").Select()
